$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (70) down into
# the new row (71) so the new row picks up the same per-column styles
# (bold/border index style in col A, date number-format in col E, plain
# default style everywhere else) without touching any values.
$ws.Range("A70:V70").Copy()
$ws.Range("A71:V71").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New match row: Araz vs Sabail, 2023-12-02
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = "azerbaijan"
$ws.Range("C71").Value = "premier-league"
$ws.Range("D71").Value = "2023-2024"
$ws.Range("E71").Value = 45262.54166666666
$ws.Range("F71").Value = "Araz"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Sabail"
$ws.Range("I71").Value = 1
$ws.Range("J71").Value = 1.81
$ws.Range("K71").Value = "01/12/2023 01:13"
$ws.Range("L71").Value = 1.63
$ws.Range("M71").Value = "02/12/2023 12:57"
$ws.Range("N71").Value = 3.22
$ws.Range("O71").Value = "01/12/2023 01:13"
$ws.Range("P71").Value = 3.67
$ws.Range("Q71").Value = "02/12/2023 12:57"
$ws.Range("R71").Value = 4.06
$ws.Range("S71").Value = "01/12/2023 01:13"
$ws.Range("T71").Value = 5.46
$ws.Range("U71").Value = "02/12/2023 12:57"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/araz-pfk-sabail/bajmtRZ9/"
